$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing data row (205) down
# across the newly appended rows (206-217), then fill in the values.
$ws.Range("A205:B205").Copy()
$ws.Range("A206:B217").PasteSpecial(-4122)

$newData = @(
    @(206, 204, "1.090707385520417E-16"),
    @(207, 205, "1.211152390500171E-16"),
    @(208, 206, "1.457167719820518E-16"),
    @(209, 207, "-1.079383496163347E-16"),
    @(210, 208, "1.561251128379126E-17"),
    @(211, 209, "-3.667701063493821E-17"),
    @(212, 210, "-5.782411586589357E-19"),
    @(213, 211, "-2.498001805406602E-17"),
    @(214, 212, "-9.562663161322149E-17"),
    @(215, 213, "-1.273576151946306E-16"),
    @(216, 214, "0"),
    @(217, 215, "0")
)

foreach ($item in $newData) {
    $rowNum = $item[0]
    $aVal = $item[1]
    $bVal = $item[2]
    $ws.Cells.Item($rowNum, 1).Value = $aVal
    $ws.Cells.Item($rowNum, 2).Value = [double]$bVal
}
